$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds a date value that was bumped from serial 45182 (2023-09-13)
# to serial 45184 (2023-09-15) for every data row (rows 2 through 238).
$ws.Range("C2:C238").Value = 45184
